$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 7.594773609502983
$ws.Range("C2").Value2 = 5.602546635525333
$ws.Range("E2").Value2 = 22.42759765640945
$ws.Range("F2").Value2 = 38.62044359518265
$ws.Range("G2").Value2 = 25.05905047559544
$ws.Range("H2").Value2 = 13.5361465823983
$ws.Range("I2").Value2 = 19.58758553010006
$ws.Range("J2").Value2 = 7.894537636169227
$ws.Range("K2").Value2 = 8.530476119556415
$ws.Range("M2").Value2 = 17.66037444595479
$ws.Range("O2").Value2 = 20.07952932358564
$ws.Range("B3").Value2 = 7.218914571698801
$ws.Range("C3").Value2 = 5.488478851139501
$ws.Range("E3").Value2 = 22.31804231308687
$ws.Range("F3").Value2 = 38.59631939499504
$ws.Range("G3").Value2 = 25.21443010011415
$ws.Range("H3").Value2 = 13.58954845042495
$ws.Range("I3").Value2 = 19.69309559922957
$ws.Range("J3").Value2 = 7.905682409279717
$ws.Range("K3").Value2 = 8.289018222716368
$ws.Range("M3").Value2 = 17.50887057974801
$ws.Range("O3").Value2 = 20.17891103143022
$ws.Range("B4").Value2 = 6.996249169979955
$ws.Range("C4").Value2 = 5.416817946187014
$ws.Range("E4").Value2 = 22.2552244398679
$ws.Range("F4").Value2 = 38.59146542948385
$ws.Range("G4").Value2 = 25.3179494749857
$ws.Range("H4").Value2 = 13.62430713073114
$ws.Range("I4").Value2 = 19.7614876931274
$ws.Range("J4").Value2 = 7.912911993165743
$ws.Range("K4").Value2 = 8.135904593211359
$ws.Range("M4").Value2 = 17.41785916199772
$ws.Range("O4").Value2 = 20.24397946105784
$ws.Range("B5").Value2 = 6.9034389580335
$ws.Range("C5").Value2 = 5.387234706936969
$ws.Range("E5").Value2 = 22.23076418087276
$ws.Range("F5").Value2 = 38.591994240939
$ws.Range("G5").Value2 = 25.36216683473507
$ws.Range("H5").Value2 = 13.63896752765086
$ws.Range("I5").Value2 = 19.79026673979159
$ws.Range("J5").Value2 = 7.915955593253583
$ws.Range("K5").Value2 = 8.072348683807522
$ws.Range("M5").Value2 = 17.38131027644443
$ws.Range("O5").Value2 = 20.27151296389938
$ws.Range("B6").Value2 = 6.887906274113076
$ws.Range("C6").Value2 = 5.382300251651582
$ws.Range("E6").Value2 = 22.22677188146927
$ws.Range("F6").Value2 = 38.59223349950204
$ws.Range("G6").Value2 = 25.36963160004101
$ws.Range("H6").Value2 = 13.64143184908957
$ws.Range("I6").Value2 = 19.79510040228066
$ws.Range("J6").Value2 = 7.916466876298786
$ws.Range("K6").Value2 = 8.061727029148576
$ws.Range("M6").Value2 = 17.37527485998144
$ws.Range("O6").Value2 = 20.27614633518924
$ws.Range("B7").Value2 = 6.995005732060852
$ws.Range("C7").Value2 = 5.416420481359596
$ws.Range("E7").Value2 = 22.25488992560476
$ws.Range("F7").Value2 = 38.59146240893893
$ws.Range("G7").Value2 = 25.3185375878033
$ws.Range("H7").Value2 = 13.62450283700764
$ws.Range("I7").Value2 = 19.76187213628571
$ws.Range("J7").Value2 = 7.912952645158825
$ws.Range("K7").Value2 = 8.135052072550444
$ws.Range("M7").Value2 = 17.41736402574745
$ws.Range("O7").Value2 = 20.24434666759866
$ws.Range("B8").Value2 = 7.462324814130389
$ws.Range("C8").Value2 = 5.563568815925613
$ws.Range("E8").Value2 = 22.38891173574743
$ws.Range("F8").Value2 = 38.61006089374312
$ws.Range("G8").Value2 = 25.11093592411933
$ws.Range("H8").Value2 = 13.55415112279882
$ws.Range("I8").Value2 = 19.62321726006273
$ws.Range("J8").Value2 = 7.898300261943528
$ws.Range("K8").Value2 = 8.448264419146524
$ws.Range("M8").Value2 = 17.60773823555103
$ws.Range("O8").Value2 = 20.11295590260344
$ws.Range("B9").Value2 = 8.368074155527266
$ws.Range("C9").Value2 = 5.838121765437099
$ws.Range("E9").Value2 = 22.68609888210649
$ws.Range("F9").Value2 = 38.72534833190472
$ws.Range("G9").Value2 = 24.76859894942467
$ws.Range("H9").Value2 = 13.43179072952696
$ws.Range("I9").Value2 = 19.37988609310141
$ws.Range("J9").Value2 = 7.872623003976971
$ws.Range("K9").Value2 = 9.021522394556214
$ws.Range("M9").Value2 = 17.9955442734629
$ws.Range("O9").Value2 = 19.88742560650606
$ws.Range("B10").Value2 = 8.968994313442431
$ws.Range("C10").Value2 = 6.029888147105962
$ws.Range("E10").Value2 = 22.92404823457025
$ws.Range("F10").Value2 = 38.85771355008528
$ws.Range("G10").Value2 = 24.55708336931445
$ws.Range("H10").Value2 = 13.35136060265258
$ws.Range("I10").Value2 = 19.21844218284438
$ws.Range("J10").Value2 = 7.855604373484957
$ws.Range("K10").Value2 = 9.414850705832615
$ws.Range("M10").Value2 = 18.28726027095177
$ws.Range("O10").Value2 = 19.74132797419729
$ws.Range("B11").Value2 = 9.22802594857529
$ws.Range("C11").Value2 = 6.114698418609207
$ws.Range("E11").Value2 = 23.03624065245907
$ws.Range("F11").Value2 = 38.92815711835993
$ws.Range("G11").Value2 = 24.46965591271798
$ws.Range("H11").Value2 = 13.3168181955866
$ws.Range("I11").Value2 = 19.14874290209453
$ws.Range("J11").Value2 = 7.848259577002289
$ws.Range("K11").Value2 = 9.587224129739877
$ws.Range("M11").Value2 = 18.42100377465212
$ws.Range("O11").Value2 = 19.67912430822642
$ws.Range("B12").Value2 = 9.324034034912824
$ws.Range("C12").Value2 = 6.146444604938521
$ws.Range("E12").Value2 = 23.079264507856
$ws.Range("F12").Value2 = 38.95628971968278
$ws.Range("G12").Value2 = 24.43782298496649
$ws.Range("H12").Value2 = 13.30403146667026
$ws.Range("I12").Value2 = 19.12288656624877
$ws.Range("J12").Value2 = 7.845535125533148
$ws.Range("K12").Value2 = 9.651517746965215
$ws.Range("M12").Value2 = 18.47176128638434
$ws.Range("O12").Value2 = 19.65618209756795
$ws.Range("B13").Value2 = 9.303449835467315
$ws.Range("C13").Value2 = 6.139624251281804
$ws.Range("E13").Value2 = 23.06997504481071
$ws.Range("F13").Value2 = 38.95016629549443
$ws.Range("G13").Value2 = 24.44462196225364
$ws.Range("H13").Value2 = 13.3067722604396
$ws.Range("I13").Value2 = 19.12843130908646
$ws.Range("J13").Value2 = 7.846119359581713
$ws.Range("K13").Value2 = 9.637715135283509
$ws.Range("M13").Value2 = 18.46082537939
$ws.Range("O13").Value2 = 19.66109583737756
$ws.Range("B14").Value2 = 9.235966411847366
$ws.Range("C14").Value2 = 6.11731771851834
$ws.Range("E14").Value2 = 23.03976961594649
$ws.Range("F14").Value2 = 38.93044248747621
$ws.Range("G14").Value2 = 24.46701141132653
$ws.Range("H14").Value2 = 13.31576033930451
$ws.Range("I14").Value2 = 19.14660492296904
$ws.Range("J14").Value2 = 7.8480342964414
$ws.Range("K14").Value2 = 9.59253345775867
$ws.Range("M14").Value2 = 18.42517760247259
$ws.Range("O14").Value2 = 19.67722454899845
$ws.Range("B15").Value2 = 9.194359247426835
$ws.Range("C15").Value2 = 6.103605593256697
$ws.Range("E15").Value2 = 23.02133724682988
$ws.Range("F15").Value2 = 38.91855040375157
$ws.Range("G15").Value2 = 24.48089178845801
$ws.Range("H15").Value2 = 13.32130403707739
$ws.Range("I15").Value2 = 19.15780672927851
$ws.Range("J15").Value2 = 7.849214648972366
$ws.Range("K15").Value2 = 9.564729630923775
$ws.Range("M15").Value2 = 18.40335574126955
$ws.Range("O15").Value2 = 19.68718370270821
$ws.Range("B16").Value2 = 8.951775725543754
$ws.Range("C16").Value2 = 6.024294975667384
$ws.Range("E16").Value2 = 22.91679320992696
$ws.Range("F16").Value2 = 38.85331460608744
$ws.Range("G16").Value2 = 24.56297470404293
$ws.Range("H16").Value2 = 13.35365915668825
$ws.Range("I16").Value2 = 19.22307242451296
$ws.Range("J16").Value2 = 7.856092343119636
$ws.Range("K16").Value2 = 9.403450475367684
$ws.Range("M16").Value2 = 18.27853751568466
$ws.Range("O16").Value2 = 19.74547885939008
$ws.Range("B17").Value2 = 8.799271806434271
$ws.Range("C17").Value2 = 5.975003990763542
$ws.Range("E17").Value2 = 22.85364991604744
$ws.Range("F17").Value2 = 38.81590542316158
$ws.Range("G17").Value2 = 24.6155886457056
$ws.Range("H17").Value2 = 13.37403160338387
$ws.Range("I17").Value2 = 19.26406868364441
$ws.Range("J17").Value2 = 7.860413116316048
$ws.Range("K17").Value2 = 9.302804018237381
$ws.Range("M17").Value2 = 18.20220509236995
$ws.Range("O17").Value2 = 19.7823319149596
$ws.Range("B18").Value2 = 8.71020767850141
$ws.Range("C18").Value2 = 5.946425999456067
$ws.Range("E18").Value2 = 22.81770495861119
$ws.Range("F18").Value2 = 38.79535250055734
$ws.Range("G18").Value2 = 24.64667776340059
$ws.Range("H18").Value2 = 13.38594186779433
$ws.Range("I18").Value2 = 19.28800094606447
$ws.Range("J18").Value2 = 7.86293569814977
$ws.Range("K18").Value2 = 9.244299988401124
$ws.Range("M18").Value2 = 18.15840089709386
$ws.Range("O18").Value2 = 19.80392944223079
$ws.Range("B19").Value2 = 8.679821312665172
$ws.Range("C19").Value2 = 5.936711635312024
$ws.Range("E19").Value2 = 22.80559960306823
$ws.Range("F19").Value2 = 38.78855956805618
$ws.Range("G19").Value2 = 24.65734574280392
$ws.Range("H19").Value2 = 13.39000756660071
$ws.Range("I19").Value2 = 19.29616453406897
$ws.Range("J19").Value2 = 7.863796229489564
$ws.Range("K19").Value2 = 9.22438715497389
$ws.Range("M19").Value2 = 18.1435879130847
$ws.Range("O19").Value2 = 19.8113107769959
$ws.Range("B20").Value2 = 8.815645779876183
$ws.Range("C20").Value2 = 5.980274762383481
$ws.Range("E20").Value2 = 22.86033318730681
$ws.Range("F20").Value2 = 38.81978803421748
$ws.Range("G20").Value2 = 24.60990214345015
$ws.Range("H20").Value2 = 13.37184299437668
$ws.Range("I20").Value2 = 19.25966810680963
$ws.Range("J20").Value2 = 7.859949294879271
$ws.Range("K20").Value2 = 9.313581928655189
$ws.Range("M20").Value2 = 18.21032071386303
$ws.Range("O20").Value2 = 19.77836737527851
$ws.Range("B21").Value2 = 9.255844564887578
$ws.Range("C21").Value2 = 6.123879878326207
$ws.Range("E21").Value2 = 23.04862728719001
$ws.Range("F21").Value2 = 38.93619642178712
$ws.Range("G21").Value2 = 24.46040043567416
$ws.Range("H21").Value2 = 13.31311235424259
$ws.Range("I21").Value2 = 19.14125231746783
$ws.Range("J21").Value2 = 7.847470291658666
$ws.Range("K21").Value2 = 9.605831304237999
$ws.Range("M21").Value2 = 18.43564548691619
$ws.Range("O21").Value2 = 19.67247051275484
$ws.Range("B22").Value2 = 9.531400524001347
$ws.Range("C22").Value2 = 6.21557104663683
$ws.Range("E22").Value2 = 23.17481623936514
$ws.Range("F22").Value2 = 39.02076161234474
$ws.Range("G22").Value2 = 24.37012247708649
$ws.Range("H22").Value2 = 13.27644038808979
$ws.Range("I22").Value2 = 19.0669919961428
$ws.Range("J22").Value2 = 7.839645879002993
$ws.Range("K22").Value2 = 9.791103456077757
$ws.Range("M22").Value2 = 18.58354217920546
$ws.Range("O22").Value2 = 19.60683406084135
$ws.Range("B23").Value2 = 9.385447557255524
$ws.Range("C23").Value2 = 6.166838177437064
$ws.Range("E23").Value2 = 23.10719034278075
$ws.Range("F23").Value2 = 38.97485624063095
$ws.Range("G23").Value2 = 24.4176224804531
$ws.Range("H23").Value2 = 13.29585640519746
$ws.Range("I23").Value2 = 19.10633989595254
$ws.Range("J23").Value2 = 7.84379167459342
$ws.Range("K23").Value2 = 9.692756030137707
$ws.Range("M23").Value2 = 18.50456139335996
$ws.Range("O23").Value2 = 19.64153822680873
$ws.Range("B24").Value2 = 8.808247426910929
$ws.Range("C24").Value2 = 5.97789259233352
$ws.Range("E24").Value2 = 22.85731056627283
$ws.Range("F24").Value2 = 38.81802973271061
$ws.Range("G24").Value2 = 24.61247039323921
$ws.Range("H24").Value2 = 13.37283184810664
$ws.Range("I24").Value2 = 19.26165647690283
$ws.Range("J24").Value2 = 7.860158868570645
$ws.Range("K24").Value2 = 9.308711229146359
$ws.Range("M24").Value2 = 18.20665138826028
$ws.Range("O24").Value2 = 19.78015846592855
$ws.Range("B25").Value2 = 8.134243374130893
$ws.Range("C25").Value2 = 5.765498443953997
$ws.Range("E25").Value2 = 22.60215593791549
$ws.Range("F25").Value2 = 38.68575416976184
$ws.Range("G25").Value2 = 24.85422319978396
$ws.Range("H25").Value2 = 13.46322683617695
$ws.Range("I25").Value2 = 19.44266398784188
$ws.Range("J25").Value2 = 7.879243936304917
$ws.Range("K25").Value2 = 9.414850705832615
$ws.Range("M25").Value2 = 17.88929003558171
$ws.Range("O25").Value2 = 19.94499688080817
